# Apply the edit described by the diff:
# - Rows 2-11: column A becomes "test_10".."test_19" (was a mix of text labels),
#   column B becomes numeric 1 (was a shared string reference).
# - Rows 12-15 are removed entirely (data shrinks from 15 rows to 11 rows).
# - Row 1 (header: title/label) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 12-15 first (from the bottom up) so row indices for rows 2-11 stay stable.
$ws.Range("A15:B15").Delete()
$ws.Range("A14:B14").Delete()
$ws.Range("A13:B13").Delete()
$ws.Range("A12:B12").Delete()

# Update rows 2 through 11: A = test_10 .. test_19, B = 1 (numeric)
for ($r = 2; $r -le 11; $r++) {
    $n = $r + 8
    $ws.Cells.Item($r, 1).Value = "test_$n"
    $ws.Cells.Item($r, 2).Value = 1
}
